$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: DF - PIB 2022 -> 2023
$ws.Range("B2").Value = 118174.1116095417
$ws.Range("D2").Value = "PIB 2023 Deflacionado"

# Row 3: MT
$ws.Range("B3").Value = 76532.28963539573
$ws.Range("D3").Value = "PIB 2023 Deflacionado"

# Row 4: SP
$ws.Range("B4").Value = 73845.19036585005
$ws.Range("D4").Value = "PIB 2023 Deflacionado"

# Row 5: was RJ -> now SC
$ws.Range("A5").Value = "SC"
$ws.Range("B5").Value = 69959.10090505773
$ws.Range("D5").Value = "PIB 2023 Deflacionado"

# Row 6: was SC -> now RJ
$ws.Range("A6").Value = "RJ"
$ws.Range("B6").Value = 67161.88535005336
$ws.Range("D6").Value = "PIB 2023 Deflacionado"

# Row 7: MS
$ws.Range("B7").Value = 64948.89321994126
$ws.Range("D7").Value = "PIB 2023 Deflacionado"

# Row 8: SE, rank 22 -> 23
$ws.Range("B8").Value = 26006.98661973922
$ws.Range("C8").Value = 23
$ws.Range("D8").Value = "PIB 2023 Deflacionado"

# Row 9: BR
$ws.Range("B9").Value = 51300.70579350938
$ws.Range("D9").Value = "PIB 2023 Deflacionado"

# Row 10: NE
$ws.Range("B10").Value = 26237.41536180414
$ws.Range("D10").Value = "PIB 2023 Deflacionado"

# Row 11: was MT -> now PI
$ws.Range("A11").Value = "PI"
$ws.Range("B11").Value = 1.537570136346218
$ws.Range("D11").Value = "Variação (%) 2023/2010"

# Row 12: was PI -> now AL
$ws.Range("A12").Value = "AL"
$ws.Range("B12").Value = 1.47614098883597
$ws.Range("D12").Value = "Variação (%) 2023/2010"

# Row 13: BA
$ws.Range("B13").Value = 1.456636142415122
$ws.Range("D13").Value = "Variação (%) 2023/2010"

# Row 14: was MS -> now PR
$ws.Range("A14").Value = "PR"
$ws.Range("B14").Value = 1.361548723808843
$ws.Range("D14").Value = "Variação (%) 2023/2010"

# Row 15: was PR -> now MT
$ws.Range("A15").Value = "MT"
$ws.Range("B15").Value = 1.347864888435147
$ws.Range("D15").Value = "Variação (%) 2023/2010"

# Row 16: was RJ -> now RS
$ws.Range("A16").Value = "RS"
$ws.Range("B16").Value = 1.330588661237085
$ws.Range("D16").Value = "Variação (%) 2023/2010"

# Row 17: SE, rank 22 -> 23
$ws.Range("B17").Value = 0.91785810919599
$ws.Range("C17").Value = 23
$ws.Range("D17").Value = "Variação (%) 2023/2010"

# Row 18: BR
$ws.Range("B18").Value = 1.205735709293767
$ws.Range("D18").Value = "Variação (%) 2023/2010"

# Row 19: NE
$ws.Range("B19").Value = 1.297733510014661
$ws.Range("D19").Value = "Variação (%) 2023/2010"
